# House Robber 2 and Search in Rotated Sorted Array updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blind 75")

# Row 9: Search in Rotated Sorted Array
$ws.Range("D9").Value = "Medium"
$ws.Range("E9").Value = "X"
$ws.Range("G9").Value = "O(log(n))"
$ws.Range("F9").Value = "While l<=r: if mid = target return. If l<mid then check if target is in range of l and mid and change r. Else change l. Else check if target inbetween mid and r and change l else change r."

# Row 25: House Robber II
$ws.Range("E25").Value = "X"
$ws.Range("F25").Value = "Do House Robber but take max between houses - first house and houses - last house"
$ws.Range("G25").Value = "O(N)"

# Update selection to match the final state
$ws.Range("F13").Select()
